# CMP73010.docx edit:
#  1. Collapse the ">> > your  stuff after this line >>>" run/proofErr
#     soup into a single plain run reading ">>>  your stuff after this
#     line >>>".
#  2. Replace the "Ben changing things up!" paragraph text with the new
#     "Version control is like ..." sentence.
#  3. Move the "_GoBack" bookmark (Word's "last edit" marker) from its
#     old spot - right after "MP73010" in the title line - to the end
#     of the paragraph that now holds the new "Version control ..."
#     text (collapsed, right after the run, before the paragraph mark).

$d = $word.ActiveDocument

# --- 1. Tidy up the ">>>  your stuff after this line >>>" paragraph ---
# The original paragraph is split across three runs with a couple of
# <w:proofErr/> markers sitting between them (an artefact of Word's
# grammar checker). A single Find/Replace across the whole phrase
# collapses it back down to one plain run and drops the proofErr
# markers in the process.
$rGt = $d.Content
$rGt.Find.Execute(
    ">>>  your stuff after this line >>>", $true, $false, $false,
    $false, $false, $true, 1, $false,
    ">>>  your stuff after this line >>>", 2) | Out-Null

# --- 2. Swap in the new "Version control ..." sentence ---
$newText = "Version control is like you are saving your assignment at each task or at each paragraph , in different files with different file names so it makes sense."
$rBen = $d.Content
$rBen.Find.Execute(
    "Ben changing things up!", $true, $false, $false,
    $false, $false, $true, 1, $false,
    $newText, 2) | Out-Null

# --- 3. Move the _GoBack bookmark ---
# 3a. Delete it from its old home (right after "MP73010").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3b. Re-create it, collapsed, immediately after the new paragraph's
# text and before its paragraph mark. A collapsed range placed exactly
# on a paragraph's own text-end boundary lands in the wrong spot, so
# work around that: temporarily type a placeholder character after the
# text, drop the (now safely mid-paragraph) bookmark right in front of
# it, then delete the placeholder again. The bookmark stays put because
# it sits to the left of the deleted character.
$targetPara = $d.Paragraphs.Item(5)
$paraEnd = $targetPara.Range.End - 1   # position just before the pilcrow

$placeholder = $d.Range($paraEnd, $paraEnd)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($paraEnd, $paraEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($paraEnd, $paraEnd + 1).Delete() | Out-Null
